# Auto-generated Excel COM-interop edit script
# Updates the crypto price/volume table to reflect the latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''27.001.74'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.54%  '

# Row 3
$ws.Range("D3").Value = '''1.823.05'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.55%  '

# Row 4
$ws.Range("E4").Value = '  +0.24%  '

# Row 5
$ws.Range("D5").Value = '''311.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.31%  '

# Row 6
$ws.Range("E6").Value = '  +0.18%  '

# Row 7
$ws.Range("D7").Value = '''0.4695'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.17%  '

# Row 8
$ws.Range("D8").Value = '''0.3664'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.86%  '

# Row 9
$ws.Range("D9").Value = '''0.07358'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.10%  '

# Row 10
$ws.Range("D10").Value = '''0.8743'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.61%  '

# Row 11
$ws.Range("E11").Value = '  -0.54%  '

# Row 12
$ws.Range("D12").Value = '''1.845.40'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.67%  '

# Row 13
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '''0.07294'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.00%  '

# Row 14
$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").Value = '''5.433'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.79%  '

# Row 15
$ws.Range("D15").Value = '''6.523'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.25%  '

# Row 16
$ws.Range("D16").Value = '''91.77'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.24%  '

# Row 17
$ws.Range("D17").Value = '''1.005'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.30%  '

# Row 18
$ws.Range("D18").Value = '''0.000008748'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.38%  '

# Row 19
$ws.Range("D19").Value = '''1.004'
$ws.Range("D19").Style = "Normal"

# Row 20
$ws.Range("E20").Value = '  +0.07%  '

# Row 21
$ws.Range("D21").Value = '''27.016.79'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.46%  '

# Row 22
$ws.Range("D22").Value = '''5.290'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.90%  '

# Row 23
$ws.Range("E23").Value = '  +0.74%  '

# Row 24
$ws.Range("D24").Value = '''2.067.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.96%  '

# Row 25
$ws.Range("D25").Value = '''1.892'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.39%  '

# Row 26
$ws.Range("D26").Value = '''150.89'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.67%  '

# Row 27
$ws.Range("D27").Value = '''18.44'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.45%  '

# Row 28
$ws.Range("D28").Value = '''2.145'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.93%  '

# Row 29
$ws.Range("D29").Value = '''5.258'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.72%  '

# Row 30
$ws.Range("D30").Value = '''116.93'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.40%  '

# Row 31
$ws.Range("D31").Value = '''0.08886'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.22%  '

# Row 32
$ws.Range("D32").Value = '''0.7547'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.28%  '

# Row 33
$ws.Range("E33").Value = '  +1.11%  '

# Row 34
$ws.Range("D34").Value = '''4.508'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.16%  '

# Row 35
$ws.Range("D35").Value = '''2.932'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.01%  '

# Row 36
$ws.Range("D36").Value = '''1.004'
$ws.Range("D36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = '''1.096'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.01%  '

# Row 38
$ws.Range("D38").Value = '''0.05312'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.13%  '

# Row 39
$ws.Range("D39").Value = '''0.01949'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.19%  '

# Row 40
$ws.Range("D40").Value = '''2.977'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.16%  '

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '''7.212'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.09%  '

# Row 42
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = '''2.374'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.24%  '

# Row 43
$ws.Range("E43").Value = '  -0.60%  '

# Row 44
$ws.Range("D44").Value = '''0.1655'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.21%  '

# Row 45
$ws.Range("D45").Value = '''8.482'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.63%  '

# Row 46
$ws.Range("D46").Value = '''0.4890'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.96%  '

# Row 47
$ws.Range("D47").Value = '''10.51'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.60%  '

# Row 48
$ws.Range("E48").Value = '  +0.19%  '

# Row 49
$ws.Range("D49").Value = '''1.663'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.38%  '

# Row 50
$ws.Range("D50").Value = '''103.21'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.32%  '

# Row 51
$ws.Range("D51").Value = '''0.06301'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.56%  '

